# "Generate Report for Archive"
#
# 1) Localization status moves from "Ready for handoff" to "In Translation"
#    for the tracked file — this shows up on every sheet that surfaces the
#    Status column (Overview!E2:F2, zh-cn!C2, de-de!C2), all of which share
#    the same underlying string.
# 2) The two now-narrower "handoff status" columns get re-sized down to fit
#    the new, shorter label (Overview E:F, and column C on each language
#    sheet).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$zhcn.Range("C2").Value = "In Translation"
$dede.Range("C2").Value = "In Translation"

# --- Column widths follow the shorter text ---
$overview.Range("E1:F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
